# Update countries & provincias Spain
# This refreshes the "Pais" (country) COVID stats sheet:
#   - updates the "last refreshed" timestamp in A1
#   - updates several countries' case numbers
#   - a handful of countries moved rows because the sheet is kept sorted
#     by "Casos totales" (column B) descending, and their updated totals
#     changed their rank order relative to their neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: last-updated timestamp -----------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 20 de Abril de 2020 a las 16:22"

# --- Simple in-place stat refreshes (country stays on the same row) ---
# Noruega (row 38)
$ws.Cells.Item(38, 5).Value = 6902   # Casos activos
$ws.Cells.Item(38, 7).Value = 4      # Muertes hoy
$ws.Cells.Item(38, 8).Value = 169    # Muertes

# Bulgaria (row 85)
$ws.Cells.Item(85, 2).Value = 929    # Casos totales
$ws.Cells.Item(85, 3).Value = 35     # Nuevos casos
$ws.Cells.Item(85, 5).Value = 719    # Recuperados
$ws.Cells.Item(85, 6).Value = 34     # Casos criticos

# Tunez (row 86)
$ws.Cells.Item(86, 4).Value = 148    # Casos activos
$ws.Cells.Item(86, 5).Value = 693    # Recuperados

# --- Rows that swapped order because of updated "Casos totales" -------

# Azerbaiyan overtakes Oman (rows 72/73)
$ws.Cells.Item(72, 1).Value = "Azerbaiyan"
$ws.Cells.Item(72, 2).Value = 1436
$ws.Cells.Item(72, 3).Value = 38
$ws.Cells.Item(72, 4).Value = 791
$ws.Cells.Item(72, 5).Value = 626
$ws.Cells.Item(72, 6).Value = 16
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 19

$ws.Cells.Item(73, 1).Value = "Oman"
$ws.Cells.Item(73, 2).Value = 1410
$ws.Cells.Item(73, 3).Value = 144
$ws.Cells.Item(73, 4).Value = 238
$ws.Cells.Item(73, 5).Value = 1165
$ws.Cells.Item(73, 6).Value = 3
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 7

# Zambia overtakes Islas Caimanes (rows 155/156)
$ws.Cells.Item(155, 1).Value = "Zambia"
$ws.Cells.Item(155, 2).Value = 65
$ws.Cells.Item(155, 3).Value = 4
$ws.Cells.Item(155, 4).Value = 35
$ws.Cells.Item(155, 5).Value = 27
$ws.Cells.Item(155, 6).Value = 1
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 3

$ws.Cells.Item(156, 1).Value = "Islas Caimanes"
$ws.Cells.Item(156, 2).Value = 61
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 7
$ws.Cells.Item(156, 5).Value = 53
$ws.Cells.Item(156, 6).Value = 3
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 1

# Siria overtakes Mozambique (rows 168/169)
$ws.Cells.Item(168, 1).Value = "Siria"
$ws.Cells.Item(168, 2).Value = 39
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 5
$ws.Cells.Item(168, 5).Value = 31
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 3

$ws.Cells.Item(169, 1).Value = "Mozambique"
$ws.Cells.Item(169, 2).Value = 39
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(169, 4).Value = 8
$ws.Cells.Item(169, 5).Value = 31
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 0

# Suazilandia overtakes Antigua y Barbuda and Timor Oriental (rows 177/178/179)
$ws.Cells.Item(177, 1).Value = "Suazilandia"
$ws.Cells.Item(177, 2).Value = 24
$ws.Cells.Item(177, 3).Value = 2
$ws.Cells.Item(177, 4).Value = 8
$ws.Cells.Item(177, 5).Value = 15
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 1

$ws.Cells.Item(178, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(178, 2).Value = 23
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 4).Value = 3
$ws.Cells.Item(178, 5).Value = 17
$ws.Cells.Item(178, 6).Value = 1
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(178, 8).Value = 3

$ws.Cells.Item(179, 1).Value = "Timor Oriental"
$ws.Cells.Item(179, 2).Value = 22
$ws.Cells.Item(179, 3).Value = 3
$ws.Cells.Item(179, 4).Value = 1
$ws.Cells.Item(179, 5).Value = 21
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 0
